$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-06-27 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-28 Saturday", 2)

# Update unique multiplication expressions (non-duplicated text)
$d.Content.Find.Execute("18×91=", $true, $false, $false, $false, $false, $true, 1, $false, "41×19=", 2)
$d.Content.Find.Execute("62×70=", $true, $false, $false, $false, $false, $true, 1, $false, "70×20=", 2)
$d.Content.Find.Execute("23×23=", $true, $false, $false, $false, $false, $true, 1, $false, "11×39=", 2)
$d.Content.Find.Execute("98×90=", $true, $false, $false, $false, $false, $true, 1, $false, "30×23=", 2)
$d.Content.Find.Execute("79×83=", $true, $false, $false, $false, $false, $true, 1, $false, "57×30=", 2)
$d.Content.Find.Execute("98×98=", $true, $false, $false, $false, $false, $true, 1, $false, "57×64=", 2)
$d.Content.Find.Execute("53×12=", $true, $false, $false, $false, $false, $true, 1, $false, "69×87=", 2)
$d.Content.Find.Execute("62×48=", $true, $false, $false, $false, $false, $true, 1, $false, "43×11=", 2)
$d.Content.Find.Execute("43×73=", $true, $false, $false, $false, $false, $true, 1, $false, "89×99=", 2)
$d.Content.Find.Execute("57×27=", $true, $false, $false, $false, $false, $true, 1, $false, "28×16=", 2)
$d.Content.Find.Execute("24×70=", $true, $false, $false, $false, $false, $true, 1, $false, "64×15=", 2)
$d.Content.Find.Execute("86×86=", $true, $false, $false, $false, $false, $true, 1, $false, "74×78=", 2)
$d.Content.Find.Execute("42×41=", $true, $false, $false, $false, $false, $true, 1, $false, "67×59=", 2)
$d.Content.Find.Execute("37×12=", $true, $false, $false, $false, $false, $true, 1, $false, "86×22=", 2)
$d.Content.Find.Execute("43×49=", $true, $false, $false, $false, $false, $true, 1, $false, "15×68=", 2)
$d.Content.Find.Execute("46×30=", $true, $false, $false, $false, $false, $true, 1, $false, "50×43=", 2)
$d.Content.Find.Execute("75×28=", $true, $false, $false, $false, $false, $true, 1, $false, "89×80=", 2)
$d.Content.Find.Execute("17×11=", $true, $false, $false, $false, $false, $true, 1, $false, "82×12=", 2)
$d.Content.Find.Execute("34×47=", $true, $false, $false, $false, $false, $true, 1, $false, "90×56=", 2)
$d.Content.Find.Execute("47×51=", $true, $false, $false, $false, $false, $true, 1, $false, "78×18=", 2)
$d.Content.Find.Execute("94×36=", $true, $false, $false, $false, $false, $true, 1, $false, "40×76=", 2)
$d.Content.Find.Execute("96×53=", $true, $false, $false, $false, $false, $true, 1, $false, "15×78=", 2)
$d.Content.Find.Execute("83×85=", $true, $false, $false, $false, $false, $true, 1, $false, "16×34=", 2)

# Handle the duplicated "28×91=" cells individually via the table:
# Row 1 (1-indexed), Column 4 -> 61×15=, Column 5 -> 99×79=
$tbl = $d.Tables.Item(1)
$cell4 = $tbl.Cell(1, 4)
$cell4.Range.Text = "61×15="
$cell5 = $tbl.Cell(1, 5)
$cell5.Range.Text = "99×79="
